# Insert a new weekly observation row at row 144 (pushing all existing
# rows 144-224 down to 145-225) and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(144).Insert()

$ws.Range("A144").Value = 3
$ws.Range("B144").Value = "Femacal de La Calera"
$ws.Range("C144").Value = "Coquimbo"
$ws.Range("D144").Value = 44879
$ws.Range("E144").Value = 5
$ws.Range("F144").Value = 100112026
$ws.Range("G144").Value = "Haba"
$ws.Range("H144").Value = "Sin especificar"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 90
$ws.Range("K144").Value = 8000
$ws.Range("L144").Value = 8500
$ws.Range("M144").Value = 8250
$ws.Range("N144").Value = "$/saco 25 kilos"
$ws.Range("O144").Value = "Provincia de Quillota"
$ws.Range("P144").Value = 330
$ws.Range("Q144").Value = 25
$ws.Range("R144").Value = "Hortaliza"
